$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.726.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.731.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("E10").Value = "  +4.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.375"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.217.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.579.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.737.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("E19").Value = "  -1.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.518"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "

$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0909"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("E29").Value = "  +2.97%  "

$ws.Range("E30").Value = "  +9.45%  "

$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("E35").Value = "  +2.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  +0.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.974"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "345.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0582"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("E47").Value = "  -1.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.92%  "

$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.40%  "

